$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: Replicate the cell-level formatting for the new rows (52-60)
# by copying from existing rows that already carry the exact same
# style indices. This keeps us reusing existing style ids instead of
# minting new ones.
# ------------------------------------------------------------------

# Row 52: divider/spacer row (copy whole row incl. stray G comment string)
$ws.Range("A42:G42").Copy($ws.Range("A52:G52"))
$ws.Rows(52).RowHeight = 4.5

# Row 53: new feature row (12 - UI), same shape as row 43 (11 - CICD pipeline)
$ws.Range("A43:G43").Copy($ws.Range("A53:G53"))

# Row 54: continuation row, same shape as row 44 (only C/E/F/G populated)
$ws.Range("C44").Copy($ws.Range("C54"))
$ws.Range("E44").Copy($ws.Range("E54"))
$ws.Range("F44").Copy($ws.Range("F54"))
$ws.Range("G44").Copy($ws.Range("G54"))
$ws.Range("C54").ClearContents()

# Row 55 & 56: continuation rows, same shape as row 46 (E/F/G populated)
$ws.Range("E46:G46").Copy($ws.Range("E55:G55"))
$ws.Range("E46:G46").Copy($ws.Range("E56:G56"))

# Row 57: closing "Passed" row, same shape as row 31 (E/F/G populated)
$ws.Range("E31:G31").Copy($ws.Range("E57:G57"))

# Rows 58-60: trailing blank spacer rows, same shape as E31/G31 (blank, styled)
$ws.Range("E31").Copy($ws.Range("E58"))
$ws.Range("G31").Copy($ws.Range("G58"))
$ws.Range("E31").Copy($ws.Range("E59"))
$ws.Range("G31").Copy($ws.Range("G59"))
$ws.Range("E31").Copy($ws.Range("E60"))
$ws.Range("G31").Copy($ws.Range("G60"))

# ------------------------------------------------------------------
# Step 2: Write the real content. Order matters here: this determines
# the order new entries are appended to the shared-strings table, so
# we write in the same order the original author must have used.
# ------------------------------------------------------------------

$ws.Range("A53").Value = 12

$ws.Range("B53").Value = "UI"
$ws.Range("C53").Value = "Can enter in both weekly income sources and weekly expense sources"
$ws.Range("D53").Value = "Successfully shows those sources"
$ws.Range("E53").Value = "Gender field missing"
$ws.Range("E54").Value = "Wrong datatime entered for time"

# This comment contains a leading apostrophe. A literal leading "'" in a
# Range.Value assignment is swallowed (same as Excel's "treat as text"
# quote-prefix on manual entry), so build it as a formula using CHAR(39)
# and then bake the formula down to a plain text value via copy/paste-values.
$q = [char]39
$ws.Range("E55").Formula = "=CHAR(39) & ""forms.CreateWeeklyExpenseSource object"" & CHAR(39) & "" has no attribute "" & CHAR(39) & ""income_type"" & CHAR(39)"
$ws.Range("E55").Copy()
$ws.Range("E55").PasteSpecial(-4163)

$ws.Range("G55").Value = "Need to change income to expense"
$ws.Range("G54").Value = "The form should automatically check for this."
$ws.Range("E56").Value = "Text box is lift align and looks weird. Should be centered"
$ws.Range("G56").Value = "Use in-line style and margin property set to auto"
